# Update NATMI TPM-derived computed values for Clcf1-Crlf1 ligand-receptor pairs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06373566666666666
$ws.Range("H2").Value = 0.191207
$ws.Range("I2").Value = 0.01058875298517695
$ws.Range("J2").Value = 0.01058875298517695
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1993696666666667
$ws.Range("N2").Value = 0.598109
$ws.Range("O2").Value = 0.007341795645341
$ws.Range("P2").Value = 0.007341795645341001
$ws.Range("Q2").Value = 0.01270695861811111
$ws.Range("R2").Value = 0.114362627563
$ws.Range("S2").Value = 0.00007774046055616368
$ws.Range("T2").Value = 0.00007774046055616368

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06373566666666666
$ws.Range("H3").Value = 0.191207
$ws.Range("I3").Value = 0.01058875298517695
$ws.Range("J3").Value = 0.01058875298517695
$ws.Range("O3").Value = 0.6504066282123248
$ws.Range("P3").Value = 0.6504066282123249
$ws.Range("Q3").Value = 1.125704188577333
$ws.Range("R3").Value = 10.131337697196
$ws.Range("S3").Value = 0.006886995126062131
$ws.Range("T3").Value = 0.006886995126062132

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06373566666666666
$ws.Range("H4").Value = 0.191207
$ws.Range("I4").Value = 0.01058875298517695
$ws.Range("J4").Value = 0.01058875298517695
$ws.Range("O4").Value = 0.3422515761423342
$ws.Range("P4").Value = 0.3422515761423342
$ws.Range("Q4").Value = 0.592358712379
$ws.Range("R4").Value = 5.331228411411
$ws.Range("S4").Value = 0.003624017398558659
$ws.Range("T4").Value = 0.003624017398558659

# Row 5
$ws.Range("I5").Value = 0.2961697031425515
$ws.Range("J5").Value = 0.2961697031425515
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1993696666666667
$ws.Range("N5").Value = 0.598109
$ws.Range("O5").Value = 0.007341795645341
$ws.Range("P5").Value = 0.007341795645341001
$ws.Range("Q5").Value = 0.3554163712232222
$ws.Range("R5").Value = 3.198747341009
$ws.Range("S5").Value = 0.002174417436813922
$ws.Range("T5").Value = 0.002174417436813922

# Row 6
$ws.Range("I6").Value = 0.2961697031425515
$ws.Range("J6").Value = 0.2961697031425515
$ws.Range("O6").Value = 0.6504066282123248
$ws.Range("P6").Value = 0.6504066282123249
$ws.Range("S6").Value = 0.1926307379995921
$ws.Range("T6").Value = 0.1926307379995921

# Row 7
$ws.Range("I7").Value = 0.2961697031425515
$ws.Range("J7").Value = 0.2961697031425515
$ws.Range("O7").Value = 0.3422515761423342
$ws.Range("P7").Value = 0.3422515761423342
$ws.Range("S7").Value = 0.1013645477061455
$ws.Range("T7").Value = 0.1013645477061455

# Row 8
$ws.Range("I8").Value = 0.6932415438722715
$ws.Range("J8").Value = 0.6932415438722715
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1993696666666667
$ws.Range("N8").Value = 0.598109
$ws.Range("O8").Value = 0.007341795645341
$ws.Range("P8").Value = 0.007341795645341001
$ws.Range("Q8").Value = 0.8319196436702222
$ws.Range("R8").Value = 7.487276793032
$ws.Range("S8").Value = 0.005089637747970915
$ws.Range("T8").Value = 0.005089637747970915

# Row 9
$ws.Range("I9").Value = 0.6932415438722715
$ws.Range("J9").Value = 0.6932415438722715
$ws.Range("O9").Value = 0.6504066282123248
$ws.Range("P9").Value = 0.6504066282123249
$ws.Range("S9").Value = 0.4508888950866705
$ws.Range("T9").Value = 0.4508888950866706

# Row 10
$ws.Range("I10").Value = 0.6932415438722715
$ws.Range("J10").Value = 0.6932415438722715
$ws.Range("O10").Value = 0.3422515761423342
$ws.Range("P10").Value = 0.3422515761423342
$ws.Range("S10").Value = 0.2372630110376301
$ws.Range("T10").Value = 0.2372630110376301

